$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in row 3 with the new "Project Plan" entry
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 43144
$ws.Range("C3").Value = 0.78749999999999998
$ws.Range("D3").Value = "Project Plan"
$ws.Range("E3").Value = "2700-Indu"
$ws.Range("F3").Value = "Added the Target Start and End dates"

# Widen column F slightly to fit the new description text
# (COM ColumnWidth of 32.666666666666664 rounds, via the runtime's
# pixel-grid quantization, to the closest achievable stored width of 33.5,
# which is the nearest representable value to the target 33.54296875)
$ws.Columns.Item(6).ColumnWidth = 32.666666666666664

# Update the active selection to the whole column F (as in the target file)
$ws.Range("F1:F1048576").Select()
